# Daily attendance processing - rotate "Recorded By" (column G) entries.
# For every data row, if the "Recorded By" cell contains more than one
# comma-separated entry, move the last entry to the front (the most
# recent recorder is listed first), except for the literal value
# "System, admin@admin.com" which is left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count
$firstRow = $usedRange.Row
$lastRowAbs = $firstRow + $lastRow - 1

for ($r = 2; $r -le $lastRowAbs; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null -and $val -ne "" -and $val -ne "System, admin@admin.com") {
        $parts = $val -split ", "
        if ($parts.Length -gt 1) {
            $lastPart = $parts[$parts.Length - 1]
            $restParts = $parts[0..($parts.Length - 2)]
            $newParts = @($lastPart) + $restParts
            $newVal = $newParts -join ", "
            $cell.Value2 = $newVal
        }
    }
}
